$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "46.814.42"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +4.20%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.264.09"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.12%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.88"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.08"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +6.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.561"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.63%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.57"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0782"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.16"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.76%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.605.20"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.260.64"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.55"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "46.780.70"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +4.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.794"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.84%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.68"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0925"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.85"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "64.99"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "248.36"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +4.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.80"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.72%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.88%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "42.45"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.64%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.70"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.84"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.21%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.78"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +8.77%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "145.44"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -4.70%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0774"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.19"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +8.30%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +9.60%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "16.17"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +18.40%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.71"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.50%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.81"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.48%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.99%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.94%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.95"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.788.26"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.77%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.00"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +18.98%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "71.53"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.77%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.58%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "93.74"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.76%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "FraxShare"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.80"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.68%  "
